$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3: B3 79243 -> 79244 ---
$ws.Range("B3").Value = 79244

# --- Rows 5 and 6 effectively swap their entire content (the sighting that
#     was on row 6 moves up to row 5, and vice versa), except the row that
#     ends up on row 5 has its "Taxonsorteringsordning" (column B) forced
#     to 79244 instead of inheriting the old row 6 value verbatim.

# Columns actually populated (even if only with a blank placeholder) in the
# original row 5 / row 6, taken from the workbook's current contents.
$row5Cols = @("A","B","D","E","F","G","H","I","K","L","M","N","P","Q","R","S","T","U","V","W","Y","Z","AA","AB","AC","AD","AE","AG","AT","AW","AX","AY")
$row6Cols = @("A","B","D","E","F","G","H","I","J","K","N","P","Q","R","S","T","U","V","W","Y","AA","AD","AE","AF","AG","AH","AT","AW","AX","AY")

# Snapshot current values (row 5 / row 6) before anything is overwritten.
$row5val = @{}
foreach ($col in $row5Cols) { $row5val[$col] = $ws.Range($col + "5").Value2 }
$row6val = @{}
foreach ($col in $row6Cols) { $row6val[$col] = $ws.Range($col + "6").Value2 }

# All columns touched by either row, so we know what to clear.
$allCols = @($row5Cols + $row6Cols | Select-Object -Unique)

# New row 5 = old row 6's data; blank out anything old row 6 didn't have.
foreach ($col in $allCols) {
    $target = $ws.Range($col + "5")
    if ($row6Cols -contains $col) {
        $target.Value = $row6val[$col]
    } else {
        $target.ClearContents()
    }
}
$ws.Range("B5").Value = 79244

# New row 6 = old row 5's data; blank out anything old row 5 didn't have.
foreach ($col in $allCols) {
    $target = $ws.Range($col + "6")
    if ($row5Cols -contains $col) {
        $target.Value = $row5val[$col]
    } else {
        $target.ClearContents()
    }
}

# --- Row 7: B7 79243 -> 79244 ---
$ws.Range("B7").Value = 79244
